{"js": "// Update the \"COMPETENCES TECHNIQUES\" skill list so the categories match\n// between resume_data and skillset.xlsx:\n//   - drop the \"Web : client\" and \"Autres : dess\" lines\n//   - fold the dropped language hint into \"Langages : d, r, python, matlab, c, c++\"\n//   - add a new \"Data Science : ...\" line (the list that used to be \"ML/AI\")\n//   - rename \"MLOps\" -> \"Machine Learning\"\n//   - rename \"ML/AI\" -> \"Autres\" with a refreshed keyword list\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraphs we care about by their current text so the script\n// does not depend on brittle absolute indices.\nlet webClientPara = null;\nlet autresDessPara = null;\nlet langagesPara = null;\nlet visualisationPara = null;\nlet mlopsPara = null;\nlet mlaiPara = null;\n\nfor (const p of paragraphs.items) {\n  const t = p.text.trim();\n  if (t === \"Web : client\") webClientPara = p;\n  else if (t === \"Autres : dess\") autresDessPara = p;\n  else if (t.startsWith(\"Langages :\")) langagesPara = p;\n  else if (t === \"Visualisation : tableau\") visualisationPara = p;\n  else if (t.startsWith(\"MLOps :\")) mlopsPara = p;\n  else if (t.startsWith(\"ML/AI :\")) mlaiPara = p;\n}\n\n// 1) Collapse the \"Web : client\" / \"Autres : dess\" / \"Langages : ...\" trio into a\n//    single \"Langages : d, r, python, matlab, c, c++\" paragraph. The very first of\n//    the three (\"Web : client\") directly follows an *empty* section-break\n//    paragraph (the <w:p> that only carries <w:sectPr>, switching the page into a\n//    2-column layout) -- deleting that paragraph's Range (text + mark) merges it\n//    with that preceding empty paragraph and silently destroys the section break.\n//    To keep the section break intact we reuse the \"Web : client\" paragraph in\n//    place (just overwrite its text) and delete the other two, unrelated,\n//    paragraphs instead -- neither of those sits next to a sectPr paragraph.\nif (webClientPara) {\n  webClientPara.insertText(\"Langages : d, r, python, matlab, c, c++\", \"Replace\");\n}\nif (autresDessPara) autresDessPara.delete();\nif (langagesPara) langagesPara.delete();\n\n// 2) Insert a new \"Data Science : ...\" paragraph right before \"Visualisation : tableau\".\nif (visualisationPara) {\n  visualisationPara.insertParagraph(\n    \"Data Science : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\",\n    \"Before\"\n  );\n}\n\n// 3) \"MLOps : ...\" -> \"Machine Learning : ...\"\nif (mlopsPara) {\n  mlopsPara.insertText(\n    \"Machine Learning : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\",\n    \"Replace\"\n  );\n}\n\n// 4) \"ML/AI : ...\" -> \"Autres : ...\" (new keyword list)\nif (mlaiPara) {\n  mlaiPara.insertText(\n    \"Autres : primes, cr\u00e9dit immobilier, bas\u00e9e \u00e0 marseille, mutuelle, gr\u00e2ce aux donn\u00e9es, e-mail ou sms\",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n", "ps1": "# Update the \"COMPETENCES TECHNIQUES\" skill list so the categories match\n# between resume_data and skillset.xlsx:\n#   - drop the \"Web : client\" and \"Autres : dess\" lines\n#   - fold the dropped language hint into \"Langages : d, r, python, matlab, c, c++\"\n#   - add a new \"Data Science : ...\" line (the list that used to be \"ML/AI\")\n#   - rename \"MLOps\" -> \"Machine Learning\"\n#   - rename \"ML/AI\" -> \"Autres\" with a refreshed keyword list\n$d = $word.ActiveDocument\n\nfunction Find-ParaByExactText($doc, [string]$text) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text.TrimEnd(\"`r\") -eq $text) {\n            return $p\n        }\n    }\n    return $null\n}\n\nfunction Find-ParaByPrefix($doc, [string]$prefix) {\n    foreach ($p in $doc.Paragraphs) {\n        if ($p.Range.Text.TrimEnd(\"`r\").StartsWith($prefix)) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# 1) Collapse the \"Web : client\" / \"Autres : dess\" / \"Langages : ...\" trio into a\n#    single \"Langages : d, r, python, matlab, c, c++\" paragraph. \"Web : client\"\n#    directly follows an *empty* section-break paragraph (a <w:p> that only\n#    carries <w:sectPr>, switching the page into a 2-column layout). Deleting\n#    that \"Web : client\" paragraph's Range (text + mark) in one shot merges it\n#    with that preceding empty paragraph and silently destroys the section\n#    break. To keep the section break intact we reuse the \"Web : client\"\n#    paragraph in place (just overwrite its text) and delete the other two\n#    paragraphs instead -- neither sits next to a sectPr paragraph.\n$webClient = Find-ParaByExactText $d \"Web : client\"\nif ($webClient -ne $null) {\n    $webClient.Range.Text = \"Langages : d, r, python, matlab, c, c++\"\n}\n\n$autresDess = Find-ParaByExactText $d \"Autres : dess\"\nif ($autresDess -ne $null) {\n    $autresDess.Range.Delete()\n}\n\n$langagesOld = Find-ParaByExactText $d \"Langages : r, python, matlab, c, c++\"\nif ($langagesOld -ne $null) {\n    $langagesOld.Range.Delete()\n}\n\n# 2) Insert a new \"Data Science : ...\" paragraph right before \"Visualisation : tableau\".\n$visu = Find-ParaByExactText $d \"Visualisation : tableau\"\nif ($visu -ne $null) {\n    $visu.Range.InsertParagraphBefore() | Out-Null\n}\n# Re-find \"Visualisation : tableau\" since the paragraph collection shifted;\n# the freshly inserted (still empty) paragraph sits right before it.\n$visu = Find-ParaByExactText $d \"Visualisation : tableau\"\nif ($visu -ne $null) {\n    $dataSciencePara = $visu.Previous()\n    $dataSciencePara.Range.Text = \"Data Science : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn\"\n}\n\n# 3) \"MLOps : ...\" -> \"Machine Learning : ...\"\n$mlops = Find-ParaByPrefix $d \"MLOps :\"\nif ($mlops -ne $null) {\n    $mlops.Range.Text = \"Machine Learning : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit\"\n}\n\n# 4) \"ML/AI : ...\" -> \"Autres : ...\" (new keyword list)\n$mlai = Find-ParaByPrefix $d \"ML/AI :\"\nif ($mlai -ne $null) {\n    $mlai.Range.Text = \"Autres : primes, cr\u00e9dit immobilier, bas\u00e9e \u00e0 marseille, mutuelle, gr\u00e2ce aux donn\u00e9es, e-mail ou sms\"\n}\n"}
